# "Submitted timesheet week 6"
# Add two new logged time entries (rows 7 & 8) to the "Week 6" sheet:
#   Row 7: Thu 2/15/2018, 3:00 PM - 5:30 PM, "Seeded states table, images table,
#          payment table, role table", 2.5 hours
#   Row 8: Thu 2/15/2018, 9:00 PM - 11:00 PM, "Seeded products table", 2 hours
# The Weekly Total (E20) / running Project Total (E21) formulas already in the
# sheet (and the ones that chain off it in every later week) recalculate
# automatically once the new hours are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 6")

# --- Row 7 -----------------------------------------------------------------
# Clone formatting (date / time / wrap-text styles) from the row above so the
# new cells pick up the same number formats instead of inventing new ones.
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(7, 1).Value = 41684                 # 2/15/2018 (serial, 1904 date system)
$ws.Cells.Item(7, 2).Value = 0.625                 # 3:00 PM
$ws.Cells.Item(7, 3).Value = 0.72916666666666663   # 5:30 PM
$ws.Cells.Item(7, 4).Value = "Seeded states table, images table, payment table, role table"
$ws.Cells.Item(7, 4).WrapText = $true
$ws.Cells.Item(7, 5).Value = 2.5
$ws.Rows.Item(7).RowHeight = 26

# --- Row 8 -------------------------------------------------------------------
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(8, 1).Value = 41684                 # 2/15/2018
$ws.Cells.Item(8, 2).Value = 0.875                 # 9:00 PM
$ws.Cells.Item(8, 3).Value = 0.95833333333333337   # 11:00 PM
$ws.Cells.Item(8, 4).Value = "Seeded products table"
$ws.Cells.Item(8, 4).WrapText = $true
$ws.Cells.Item(8, 5).Value = 2
$ws.Rows.Item(8).RowHeight = 18

# Update the sheet's recorded selection (matches what Excel would leave behind
# after the last edit landed on D8).
[void]$ws.Activate()
[void]$ws.Range("D8").Select()
